# Apply cryptos list update (commit: Sat Apr  1 09:10:47 UTC 2023, GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "28.488.33"
$ws.Range("E2").Value2 = "  +2.38%  "
$ws.Range("D3").Value2 = "1.828.99"
$ws.Range("E3").Value2 = "  +2.09%  "
$ws.Range("E4").Value2 = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "315.62"
$ws.Range("E5").Value2 = "  -0.18%  "
$ws.Range("E6").Value2 = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.5042"
$ws.Range("E7").Value2 = "  -5.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.3911"
$ws.Range("E8").Value2 = "  +1.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.07707"
$ws.Range("E9").Value2 = "  +3.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "41.90"
$ws.Range("E10").Value2 = "  +1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "1.113"
$ws.Range("E11").Value2 = "  +2.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "21.04"
$ws.Range("E12").Value2 = "  +3.58%  "
$ws.Range("B13").Value2 = "BinanceUSD"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "1.002"
$ws.Range("E13").Value2 = "  +0.22%  "
$ws.Range("B14").Value2 = "Polkadot"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "6.243"
$ws.Range("E14").Value2 = "  +1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "7.561"
$ws.Range("E15").Value2 = "  +1.34%  "
$ws.Range("D16").Value2 = "1.825.65"
$ws.Range("E16").Value2 = "  +2.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "93.23"
$ws.Range("E17").Value2 = "  +5.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "0.00001081"
$ws.Range("E18").Value2 = "  +2.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "0.06611"
$ws.Range("E19").Value2 = "  +1.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "17.71"
$ws.Range("E20").Value2 = "  +2.77%  "
$ws.Range("E21").Value2 = "  +0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "6.139"
$ws.Range("E22").Value2 = "  +3.01%  "
$ws.Range("D23").Value2 = "28.518.07"
$ws.Range("E23").Value2 = "  +2.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "11.15"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "2.260"
$ws.Range("E25").Value2 = "  +8.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "156.75"
$ws.Range("E26").Value2 = "  -0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "20.58"
$ws.Range("E27").Value2 = "  +2.10%  "
$ws.Range("D28").Value2 = "2.034.52"
$ws.Range("E28").Value2 = "  +2.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "2.405"
$ws.Range("E29").Value2 = "  +3.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "125.02"
$ws.Range("E30").Value2 = "  +2.96%  "
$ws.Range("E31").Value2 = "  +3.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "0.1087"
$ws.Range("E32").Value2 = "  -0.36%  "
$ws.Range("E33").Value2 = "  +2.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "3.667"
$ws.Range("E34").Value2 = "  +0.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "0.07067"
$ws.Range("E35").Value2 = "  +2.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.2216"
$ws.Range("E36").Value2 = "  +1.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "8.931"
$ws.Range("E37").Value2 = "  +6.29%  "
$ws.Range("E38").Value2 = "  +2.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "5.125"
$ws.Range("E39").Value2 = "  +1.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.6228"
$ws.Range("E40").Value2 = "  +2.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "11.19"
$ws.Range("E41").Value2 = "  -2.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "1.192"
$ws.Range("E42").Value2 = "  +2.11%  "
$ws.Range("E43").Value2 = "  +0.03%  "
$ws.Range("E44").Value2 = "  -0.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "13.47"
$ws.Range("E45").Value2 = "  +2.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "3.714"
$ws.Range("E46").Value2 = "  +1.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.5887"
$ws.Range("E47").Value2 = "  +3.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "124.26"
$ws.Range("E48").Value2 = "  -0.56%  "
$ws.Range("E49").Value2 = "  +3.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "1.186"
$ws.Range("E50").Value2 = "  +1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.06931"
$ws.Range("E51").Value2 = "  +2.07%  "
